$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "sheet_one"

# Set column widths (D and E) to match the target layout/indentation.
# (ColumnWidth is quantized to whole-pixel increments by the engine, so these
# inputs are chosen to land on the closest achievable width to 8.4 / 22.8.)
$ws.Columns.Item(4).ColumnWidth = 7.5
$ws.Columns.Item(5).ColumnWidth = 22

# Add navigation hyperlinks from the "Methods/Children" entries that are
# themselves classes back to where those classes are defined.
$ws.Hyperlinks.Add($ws.Range("E10"), "TractorPesticides")
$ws.Hyperlinks.Add($ws.Range("E16"), "Car")
$ws.Hyperlinks.Add($ws.Range("E17"), "Bike")

# Adding a hyperlink auto-applies Excel's built-in "Hyperlink" cell style
# (blue/underlined font). Restore the original (default) cell style so the
# cells keep their prior appearance - only the navigation target is new.
$ws.Range("E10").Style = "Normal"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Style = "Normal"
